$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 10699.4
$ws.Cells.Item(32, 10).Value = 10699.4
$ws.Cells.Item(32, 12).Value = 10699.4
$ws.Cells.Item(32, 14).Value = -11351.4

$ws.Cells.Item(105, 8).Value = 10000
$ws.Cells.Item(105, 10).Value = 10000
$ws.Cells.Item(105, 12).Value = 10000
$ws.Cells.Item(105, 14).Value = -16988

$ws.Cells.Item(141, 8).Value = 21699
$ws.Cells.Item(141, 9).Value = 5598.6665
$ws.Cells.Item(141, 11).Value = 16795.9995
$ws.Cells.Item(141, 13).Value = -11615.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3416.25
$ws.Cells.Item(61, 9).Value = 3454.182
$ws.Cells.Item(61, 11).Value = 3454.182
$ws.Cells.Item(61, 13).Value = -3242.182

$ws.Cells.Item(74, 8).Value = 2949.625
$ws.Cells.Item(74, 9).Value = 1291.1666
$ws.Cells.Item(74, 10).Value = 7925
$ws.Cells.Item(74, 11).Value = 1291.1666
$ws.Cells.Item(74, 12).Value = 7925
$ws.Cells.Item(74, 13).Value = -417.1666
$ws.Cells.Item(74, 14).Value = -9673

$ws.Cells.Item(77, 8).Value = 2949.625
$ws.Cells.Item(77, 9).Value = 1291.1666
$ws.Cells.Item(77, 10).Value = 7925
$ws.Cells.Item(77, 11).Value = 6455.833000000001
$ws.Cells.Item(77, 12).Value = 39625
$ws.Cells.Item(77, 13).Value = -2087.833000000001
$ws.Cells.Item(77, 14).Value = -48361

$ws.Cells.Item(106, 8).Value = 19749.5
$ws.Cells.Item(106, 10).Value = 19749.5
$ws.Cells.Item(106, 12).Value = 19749.5
$ws.Cells.Item(106, 14).Value = -22273.5

$ws.Cells.Item(122, 8).Value = 1753.3125
$ws.Cells.Item(122, 9).Value = 1673.3077
$ws.Cells.Item(122, 11).Value = 5019.9231
$ws.Cells.Item(122, 13).Value = -2569.9231

$ws.Cells.Item(136, 8).Value = 3416.25
$ws.Cells.Item(136, 9).Value = 3454.182
$ws.Cells.Item(136, 11).Value = 10362.546
$ws.Cells.Item(136, 13).Value = -7812.545999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 14).ClearContents()
$ws.Cells.Item(20, 8).Value = 2000
$ws.Cells.Item(20, 9).Value = 2000
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 2000
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = -1753

$ws.Cells.Item(94, 8).Value = 3777.7144
$ws.Cells.Item(94, 9).Value = 2543.2222
$ws.Cells.Item(94, 11).Value = 2543.2222
$ws.Cells.Item(94, 13).Value = -2092.2222

$ws.Cells.Item(134, 8).Value = 3152.3125
$ws.Cells.Item(134, 9).Value = 2174.0715
$ws.Cells.Item(134, 11).Value = 6522.2145
$ws.Cells.Item(134, 13).Value = -3987.2145

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1116.6666
$ws.Cells.Item(22, 9).Value = 675
$ws.Cells.Item(22, 11).Value = 675
$ws.Cells.Item(22, 13).Value = -325

$ws.Cells.Item(99, 8).Value = 3375
$ws.Cells.Item(99, 9).Value = 2000
$ws.Cells.Item(99, 11).Value = 2000
$ws.Cells.Item(99, 13).Value = -502

$ws.Cells.Item(126, 8).Value = 3375
$ws.Cells.Item(126, 9).Value = 2000
$ws.Cells.Item(126, 11).Value = 6000
$ws.Cells.Item(126, 13).Value = -3530

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 256.42856
$ws.Cells.Item(12, 10).Value = 23
$ws.Cells.Item(12, 12).Value = 69
$ws.Cells.Item(12, 14).Value = -415

$ws.Cells.Item(23, 8).Value = 573.7778
$ws.Cells.Item(23, 9).Value = 102.8
$ws.Cells.Item(23, 10).Value = 1162.5
$ws.Cells.Item(23, 11).Value = 308.4
$ws.Cells.Item(23, 12).Value = 3487.5
$ws.Cells.Item(23, 13).Value = -73.39999999999998
$ws.Cells.Item(23, 14).Value = -3957.5

$ws.Cells.Item(39, 14).ClearContents()
$ws.Cells.Item(39, 8).Value = 3251.2
$ws.Cells.Item(39, 9).Value = 3251.2
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 9753.599999999999
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).Value = -9459.599999999999

$ws.Cells.Item(55, 8).Value = 1768.75
$ws.Cells.Item(55, 9).Value = 1306.4286
$ws.Cells.Item(55, 10).Value = 5005
$ws.Cells.Item(55, 11).Value = 3919.2858
$ws.Cells.Item(55, 12).Value = 15015
$ws.Cells.Item(55, 13).Value = -3742.2858
$ws.Cells.Item(55, 14).Value = -15369

$ws.Cells.Item(113, 8).Value = 827
$ws.Cells.Item(113, 9).Value = 825
$ws.Cells.Item(113, 10).Value = 829
$ws.Cells.Item(113, 11).Value = 2475
$ws.Cells.Item(113, 12).Value = 2487
$ws.Cells.Item(113, 13).Value = -305
$ws.Cells.Item(113, 14).Value = -6827

$ws.Cells.Item(122, 8).Value = 587.5
$ws.Cells.Item(122, 10).Value = 500
$ws.Cells.Item(122, 12).Value = 4500
$ws.Cells.Item(122, 14).Value = -9400

$ws.Cells.Item(133, 8).Value = 14499
$ws.Cells.Item(133, 9).Value = 14499
$ws.Cells.Item(133, 11).Value = 43497
$ws.Cells.Item(133, 13).Value = -38437

$ws.Cells.Item(138, 8).Value = 3083.6365
$ws.Cells.Item(138, 9).Value = 2880
$ws.Cells.Item(138, 10).Value = 4000
$ws.Cells.Item(138, 11).Value = 8640
$ws.Cells.Item(138, 12).Value = 12000
$ws.Cells.Item(138, 13).Value = -3500
$ws.Cells.Item(138, 14).Value = -22280

$ws.Cells.Item(141, 8).Value = 9599.666999999999
$ws.Cells.Item(141, 9).Value = 1899.5
$ws.Cells.Item(141, 10).Value = 25000
$ws.Cells.Item(141, 11).Value = 5698.5
$ws.Cells.Item(141, 12).Value = 75000
$ws.Cells.Item(141, 13).Value = -518.5
$ws.Cells.Item(141, 14).Value = -85360

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 104.92308
$ws.Cells.Item(2, 9).Value = 110
$ws.Cells.Item(2, 11).Value = 110
$ws.Cells.Item(2, 13).Value = 3

$ws.Cells.Item(80, 8).Value = 2250
$ws.Cells.Item(80, 10).Value = 2500
$ws.Cells.Item(80, 12).Value = 2500
$ws.Cells.Item(80, 14).Value = -4496

$ws.Cells.Item(83, 8).Value = 2250
$ws.Cells.Item(83, 10).Value = 2500
$ws.Cells.Item(83, 12).Value = 12500
$ws.Cells.Item(83, 14).Value = -22484

$ws.Cells.Item(102, 8).Value = 1850.44
$ws.Cells.Item(102, 9).Value = 2012.5714
$ws.Cells.Item(102, 10).Value = 999.25
$ws.Cells.Item(102, 11).Value = 2012.5714
$ws.Cells.Item(102, 12).Value = 999.25
$ws.Cells.Item(102, 13).Value = -390.5714
$ws.Cells.Item(102, 14).Value = -4243.25

$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 1504
$ws.Cells.Item(122, 9).Value = 1504
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4512
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -2062

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6749.25
$ws.Cells.Item(7, 9).Value = 6999.2856
$ws.Cells.Item(7, 10).Value = 6399.2
$ws.Cells.Item(7, 11).Value = 6999.2856
$ws.Cells.Item(7, 12).Value = 6399.2
$ws.Cells.Item(7, 13).Value = -6887.2856
$ws.Cells.Item(7, 14).Value = -6623.2

$ws.Cells.Item(46, 8).Value = 4524.592
$ws.Cells.Item(46, 9).Value = 2948.5
$ws.Cells.Item(46, 10).Value = 4744.5117
$ws.Cells.Item(46, 11).Value = 2948.5
$ws.Cells.Item(46, 12).Value = 4744.5117
$ws.Cells.Item(46, 13).Value = -2760.5
$ws.Cells.Item(46, 14).Value = -5120.5117

$ws.Cells.Item(55, 8).Value = 2964.7144
$ws.Cells.Item(55, 9).Value = 2458.8333
$ws.Cells.Item(55, 11).Value = 2458.8333
$ws.Cells.Item(55, 13).Value = -2285.8333

$ws.Cells.Item(126, 8).Value = 6749.25
$ws.Cells.Item(126, 9).Value = 6999.2856
$ws.Cells.Item(126, 10).Value = 6399.2
$ws.Cells.Item(126, 11).Value = 20997.8568
$ws.Cells.Item(126, 12).Value = 19197.6
$ws.Cells.Item(126, 13).Value = -18527.8568
$ws.Cells.Item(126, 14).Value = -24137.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 1681733.4
$ws.Cells.Item(24, 10).Value = 22600
$ws.Cells.Item(24, 12).Value = 22600
$ws.Cells.Item(24, 14).Value = -23060

$ws.Cells.Item(28, 13).ClearContents()
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 11).Value = 0

$ws.Cells.Item(100, 8).Value = 591.6667
$ws.Cells.Item(100, 9).Value = 640
$ws.Cells.Item(100, 11).Value = 1280
$ws.Cells.Item(100, 13).Value = -739

$ws.Cells.Item(104, 8).Value = 33998
$ws.Cells.Item(104, 10).Value = 33998
$ws.Cells.Item(104, 12).Value = 33998
$ws.Cells.Item(104, 14).Value = -40986

$ws.Cells.Item(113, 8).Value = 899.75
$ws.Cells.Item(113, 9).Value = 933.3333
$ws.Cells.Item(113, 11).Value = 2799.9999
$ws.Cells.Item(113, 13).Value = -629.9998999999998
